# This edit inserts 3 new weekly price-report rows ("Tuna" / "Provincia de
# Talca", fecha 44627) right above the existing row 343 block, pushing all
# the subsequent rows down by three (old row 343 -> new row 346, ...,
# old row 451 -> new row 454). The inserted rows reuse the constant columns
# (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría, Unidad de
# comercialización, Kg o Unidades, Clasificación) that are common to every
# row in this subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 343:451 down to 346:454, carrying formatting (incl. the date
# number-format on column D) down with them.
$ws.Rows("343:345").Insert()

# Values shared by every row of this report subset.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$catId     = 100112027
$categoria = "Melón"
$unidadCom = "`$/unidad"
$kgOUnid   = 1
$clasif    = "Hortaliza"

$fecha    = 44627
$variedad = "Tuna"

# New rows: Extra / Primera / Segunda for "Tuna" @ "Provincia de Talca".
$newRows = @(
    @{ Row = 343; Calidad = "Extra";   Volumen = 380; PMin = 1000; PMax = 1000; PProm = 1000; Origen = "Provincia de Talca"; PKg = 1000 },
    @{ Row = 344; Calidad = "Primera"; Volumen = 420; PMin = 700;  PMax = 700;  PProm = 700;  Origen = "Provincia de Talca"; PKg = 700 },
    @{ Row = 345; Calidad = "Segunda"; Volumen = 400; PMin = 500;  PMax = 500;  PProm = 500;  Origen = "Provincia de Talca"; PKg = 500 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidadCom
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgOUnid
    $ws.Cells.Item($row, 18).Value = $clasif
}
